# Append a new order line (row 28) to the Sysco bakery order sheet.
# Every existing column in this sheet is authored as plain text (SKU
# codes, quantities and costs alike), so the new cells are written with
# a leading apostrophe to force text entry instead of letting Excel
# auto-infer numeric types for the numeric-looking values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 28

$ws.Cells.Item($newRow, 1).Value = "'7086864"
$ws.Cells.Item($newRow, 2).Value = "Chocolate Chips - White"
$ws.Cells.Item($newRow, 3).Value = "'1"
$ws.Cells.Item($newRow, 4).Value = "'99.99"
$ws.Cells.Item($newRow, 5).Value = "'99.99"
